$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 (pushes existing rows 51..115 down to 52..116)
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new data record
$ws.Range("A51").Value = 2
$ws.Range("B51").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44629
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 100112024
$ws.Range("G51").Value = "Choclo"
$ws.Range("H51").Value = "Choclero"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 50000
$ws.Range("K51").Value = 180
$ws.Range("L51").Value = 200
$ws.Range("M51").Value = 190
$ws.Range("N51").Value = "$/unidad"
$ws.Range("O51").Value = "Provincia de Limarí"
$ws.Range("P51").Value = 190
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"
